# Update the cryptos worksheet to the refreshed price/volume snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get a plain decimal value (e.g. "235.78"). Excel
# auto-converts a bare numeric-looking string typed into a General cell
# into a floating point number, which would lose the exact decimal text
# (and the text cell type) that the source data uses. Format them as
# Text first so the assigned value is kept verbatim as a string.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the refreshed values (prices, hourly volume deltas, and the
# Aave/Maker ranking swap) cell by cell, in sheet order.
$ws.Range("D2").Value = "37.387.14"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").Value = "2.063.86"
$ws.Range("E3").Value = "  +3.71%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "235.78"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("D7").Value = "58.28"
$ws.Range("E7").Value = "  +6.42%  "
$ws.Range("D9").Value = "0.383"
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("D10").Value = "58.37"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("E12").Value = "  +3.04%  "
$ws.Range("D13").Value = "2.366.89"
$ws.Range("E13").Value = "  +3.71%  "
$ws.Range("D14").Value = "14.62"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").Value = "20.98"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("D16").Value = "0.779"
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "5.20"
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").Value = "2.061.80"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").Value = "37.593.58"
$ws.Range("E19").Value = "  +3.09%  "
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  +17.43%  "
$ws.Range("D21").Value = "69.10"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").Value = "226.68"
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").Value = "164.70"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("E28").Value = "  +13.92%  "
$ws.Range("D29").Value = "8.88"
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").Value = "19.23"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("E33").Value = "  +3.11%  "
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("D36").Value = "4.51"
$ws.Range("E36").Value = "  +6.36%  "
$ws.Range("D37").Value = "3.43"
$ws.Range("E37").Value = "  +5.01%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "1.78"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "5.90"
$ws.Range("E40").Value = "  +7.99%  "
$ws.Range("D41").Value = "0.0984"
$ws.Range("E41").Value = "  +7.18%  "
$ws.Range("D42").Value = "2.96"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").Value = "4.50"
$ws.Range("E43").Value = "  +23.76%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.477.61"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "96.85"
$ws.Range("E45").Value = "  +9.00%  "
$ws.Range("D46").Value = "1.17"
$ws.Range("E46").Value = "  +5.63%  "
$ws.Range("E47").Value = "  +4.24%  "
$ws.Range("E48").Value = "  +6.56%  "
$ws.Range("E49").Value = "  +3.46%  "
$ws.Range("D50").Value = "7.30"
$ws.Range("E50").Value = "  +7.32%  "
$ws.Range("E51").Value = "  +1.68%  "
